$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-05"

# Update header label in I1 (shared string "2022 (through 11-04)" -> "2022 (through 11-05)")
$ws.Range("I1").Value = "2022 (through 11-05)"

# Update November value (row 12) from 13 to 14
$ws.Range("I12").Value = 14

# Update Total value (row 14) from 1413 to 1414
$ws.Range("I14").Value = 1414
